$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: rename "Documento" columns to the iteration/phase names ---
$ws.Range("B1").Value = "Inicio"
$ws.Range("D1").Value = "Elaboración Iteración 1"
$ws.Range("F1").Value = "Elaboración Iteración 2"

# F1 picks up the same themed border/fill style as the E column header (s=9)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- Extend the left-hand (meeting) column formatting down through row 23 ---
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B23").PasteSpecial(-4122)

# --- Add the new document row 23 (E222) under the existing E220/E221 rows ---
$ws.Range("E22:F22").Copy()
$ws.Range("E23:F23").PasteSpecial(-4122)
$ws.Range("E23").Value = "E222"

$excel.CutCopyMode = $false

# --- Update the visible selection / scroll position ---
$ws.Range("G4").Select()
